# Daily attendance processing - 2025-10-22 15:20:49
# Swap the first two comma-separated entries in the "Recorded By" column (G)
# whenever the entry begins with "System, " (e.g. "System, user@x.com" ->
# "user@x.com, System"), leaving the rest of the list untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "System, *") {
        $parts = $val -split ", "
        if ($parts.Length -ge 2 -and $parts[0] -eq "System") {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value2 = ($parts -join ", ")
        }
    }
}
